$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 12: only the "B" (catalogNumber-like) value changes
$ws.Range("B12").Value = 77685

# Row 13: values are replaced with the final (post-edit) contents
$ws.Range("A13").Value = 112182730
$ws.Range("B13").Value = 95707
$ws.Range("E13").Value = 221941
$ws.Range("F13").Value = "Plattlummer"
$ws.Range("G13").Value = "Lycopodium complanatum"
$ws.Range("H13").Value = "L."
$ws.Range("Q13").Value = 375047
$ws.Range("R13").Value = 6871264

# Row 14: values are replaced with the final (post-edit) contents
$ws.Range("A14").Value = 112182046
$ws.Range("B14").Value = 90826
$ws.Range("E14").Value = 4366
$ws.Range("F14").Value = "Skarp dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum peckii"
$ws.Range("H14").Value = "Banker"
$ws.Range("Q14").Value = 374850
$ws.Range("R14").Value = 6871061
